# Add support for text formula cells.
#
# Mirrors the upstream commit: a new row (row 7) is appended to Sheet1
# holding a formula that evaluates to a *string* result ("String Formula"),
# together with the literal expected value in column B, so the test
# fixture can exercise text-formula-cell evaluation (previously date
# formatting was checked before the formula was evaluated, throwing
# "Cannot get a numeric value from a text formula cell").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New formula cell: concatenates "String" & " " & "Formula" -> "String Formula".
$ws.Range("A7").Formula = '="String"&" "&"Formula"'

# Companion "expected value" cell, mirroring the pattern used by every
# other row in the sheet (formula in column A, expected literal in B).
$ws.Range("B7").Value = "String Formula"

# Column A now holds text wide enough to need its own explicit width
# (previously it just used the sheet default).
$ws.Columns.Item(1).ColumnWidth = 12.330729166666666

# The active selection moves on, as recorded when the fixture was saved.
$ws.Range("C12").Select()
